$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Espoo"
$ws.Range("A4").Value = "Vantaa"

$ws.Range("A5").Select()
